$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,7).Value = [double]"41.52772766666666"
$ws.Cells.Item(2,8).Value = [double]"124.583183"
$ws.Cells.Item(2,9).Value = [double]"0.750469325006714"
$ws.Cells.Item(2,10).Value = [double]"0.7504693250067139"
$ws.Cells.Item(2,13).Value = [double]"0.043414"
$ws.Cells.Item(2,15).Value = [double]"0.00150816245073376"
$ws.Cells.Item(2,16).Value = [double]"0.00150816245073376"
$ws.Cells.Item(2,17).Value = [double]"1.802884768920667"
$ws.Cells.Item(2,18).Value = [double]"16.225962920286"
$ws.Cells.Item(2,19).Value = [double]"0.001131829656402636"
$ws.Cells.Item(2,20).Value = [double]"0.001131829656402636"

$ws.Cells.Item(3,7).Value = [double]"41.52772766666666"
$ws.Cells.Item(3,8).Value = [double]"124.583183"
$ws.Cells.Item(3,9).Value = [double]"0.750469325006714"
$ws.Cells.Item(3,10).Value = [double]"0.7504693250067139"
$ws.Cells.Item(3,13).Value = [double]"7.896729"
$ws.Cells.Item(3,14).Value = [double]"23.690187"
$ws.Cells.Item(3,15).Value = [double]"0.2743251062196607"
$ws.Cells.Item(3,16).Value = [double]"0.2743251062196607"
$ws.Cells.Item(3,17).Value = [double]"327.9332113694689"
$ws.Cells.Item(3,18).Value = [double]"2951.398902325221"
$ws.Cells.Item(3,19).Value = [double]"0.2058725772970639"
$ws.Cells.Item(3,20).Value = [double]"0.2058725772970638"

$ws.Cells.Item(4,7).Value = [double]"41.52772766666666"
$ws.Cells.Item(4,8).Value = [double]"124.583183"
$ws.Cells.Item(4,9).Value = [double]"0.750469325006714"
$ws.Cells.Item(4,10).Value = [double]"0.7504693250067139"
$ws.Cells.Item(4,13).Value = [double]"2.710967666666666"
$ws.Cells.Item(4,14).Value = [double]"8.132902999999999"
$ws.Cells.Item(4,15).Value = [double]"0.09417652462385363"
$ws.Cells.Item(4,16).Value = [double]"0.09417652462385362"
$ws.Cells.Item(4,17).Value = [double]"112.5803269744721"
$ws.Cells.Item(4,18).Value = [double]"1013.222942770249"
$ws.Cells.Item(4,19).Value = [double]"0.07067659286594162"
$ws.Cells.Item(4,20).Value = [double]"0.07067659286594159"

$ws.Cells.Item(5,7).Value = [double]"41.52772766666666"
$ws.Cells.Item(5,8).Value = [double]"124.583183"
$ws.Cells.Item(5,9).Value = [double]"0.750469325006714"
$ws.Cells.Item(5,10).Value = [double]"0.7504693250067139"
$ws.Cells.Item(5,13).Value = [double]"18.134913"
$ws.Cells.Item(5,14).Value = [double]"54.404739"
$ws.Cells.Item(5,15).Value = [double]"0.629990206705752"
$ws.Cells.Item(5,16).Value = [double]"0.6299902067057519"
$ws.Cells.Item(5,17).Value = [double]"753.101728322693"
$ws.Cells.Item(5,18).Value = [double]"6777.915554904236"
$ws.Cells.Item(5,19).Value = [double]"0.472788325187306"
$ws.Cells.Item(5,20).Value = [double]"0.4727883251873058"

$ws.Cells.Item(6,9).Value = [double]"0.03901417330949282"
$ws.Cells.Item(6,10).Value = [double]"0.03901417330949282"
$ws.Cells.Item(6,13).Value = [double]"0.043414"
$ws.Cells.Item(6,15).Value = [double]"0.00150816245073376"
$ws.Cells.Item(6,16).Value = [double]"0.00150816245073376"
$ws.Cells.Item(6,17).Value = [double]"0.09372542819266667"
$ws.Cells.Item(6,18).Value = [double]"0.8435288537340001"
$ws.Cells.Item(6,19).Value = [double]"5.883971123179632E-05"
$ws.Cells.Item(6,20).Value = [double]"5.883971123179632E-05"

$ws.Cells.Item(7,9).Value = [double]"0.03901417330949282"
$ws.Cells.Item(7,10).Value = [double]"0.03901417330949282"
$ws.Cells.Item(7,13).Value = [double]"7.896729"
$ws.Cells.Item(7,14).Value = [double]"23.690187"
$ws.Cells.Item(7,15).Value = [double]"0.2743251062196607"
$ws.Cells.Item(7,16).Value = [double]"0.2743251062196607"
$ws.Cells.Item(7,17).Value = [double]"17.048056084361"
$ws.Cells.Item(7,18).Value = [double]"153.432504759249"
$ws.Cells.Item(7,19).Value = [double]"0.01070256723719887"
$ws.Cells.Item(7,20).Value = [double]"0.01070256723719887"

$ws.Cells.Item(8,9).Value = [double]"0.03901417330949282"
$ws.Cells.Item(8,10).Value = [double]"0.03901417330949282"
$ws.Cells.Item(8,13).Value = [double]"2.710967666666666"
$ws.Cells.Item(8,14).Value = [double]"8.132902999999999"
$ws.Cells.Item(8,15).Value = [double]"0.09417652462385363"
$ws.Cells.Item(8,16).Value = [double]"0.09417652462385362"
$ws.Cells.Item(8,17).Value = [double]"5.852642128686778"
$ws.Cells.Item(8,18).Value = [double]"52.673779158181"
$ws.Cells.Item(8,19).Value = [double]"0.003674219253360744"
$ws.Cells.Item(8,20).Value = [double]"0.003674219253360743"

$ws.Cells.Item(9,9).Value = [double]"0.03901417330949282"
$ws.Cells.Item(9,10).Value = [double]"0.03901417330949282"
$ws.Cells.Item(9,13).Value = [double]"18.134913"
$ws.Cells.Item(9,14).Value = [double]"54.404739"
$ws.Cells.Item(9,15).Value = [double]"0.629990206705752"
$ws.Cells.Item(9,16).Value = [double]"0.6299902067057519"
$ws.Cells.Item(9,17).Value = [double]"39.151022392817"
$ws.Cells.Item(9,18).Value = [double]"352.359201535353"
$ws.Cells.Item(9,19).Value = [double]"0.02457854710770141"
$ws.Cells.Item(9,20).Value = [double]"0.02457854710770141"

$ws.Cells.Item(10,9).Value = [double]"0.001546645353191641"
$ws.Cells.Item(10,10).Value = [double]"0.001546645353191641"
$ws.Cells.Item(10,13).Value = [double]"0.043414"
$ws.Cells.Item(10,15).Value = [double]"0.00150816245073376"
$ws.Cells.Item(10,16).Value = [double]"0.00150816245073376"
$ws.Cells.Item(10,17).Value = [double]"0.003715572718666666"
$ws.Cells.Item(10,18).Value = [double]"0.033440154468"
$ws.Cells.Item(10,19).Value = [double]"2.332592446285486E-06"
$ws.Cells.Item(10,20).Value = [double]"2.332592446285486E-06"

$ws.Cells.Item(11,9).Value = [double]"0.001546645353191641"
$ws.Cells.Item(11,10).Value = [double]"0.001546645353191641"
$ws.Cells.Item(11,13).Value = [double]"7.896729"
$ws.Cells.Item(11,14).Value = [double]"23.690187"
$ws.Cells.Item(11,15).Value = [double]"0.2743251062196607"
$ws.Cells.Item(11,16).Value = [double]"0.2743251062196607"
$ws.Cells.Item(11,17).Value = [double]"0.6758389192219999"
$ws.Cells.Item(11,18).Value = [double]"6.082550272997999"
$ws.Cells.Item(11,19).Value = [double]"0.0004242836507984416"
$ws.Cells.Item(11,20).Value = [double]"0.0004242836507984415"

$ws.Cells.Item(12,9).Value = [double]"0.001546645353191641"
$ws.Cells.Item(12,10).Value = [double]"0.001546645353191641"
$ws.Cells.Item(12,13).Value = [double]"2.710967666666666"
$ws.Cells.Item(12,14).Value = [double]"8.132902999999999"
$ws.Cells.Item(12,15).Value = [double]"0.09417652462385363"
$ws.Cells.Item(12,16).Value = [double]"0.09417652462385362"
$ws.Cells.Item(12,17).Value = [double]"0.2320172640957777"
$ws.Cells.Item(12,18).Value = [double]"2.088155376861999"
$ws.Cells.Item(12,19).Value = [double]"0.0001456576841892214"
$ws.Cells.Item(12,20).Value = [double]"0.0001456576841892213"

$ws.Cells.Item(13,9).Value = [double]"0.001546645353191641"
$ws.Cells.Item(13,10).Value = [double]"0.001546645353191641"
$ws.Cells.Item(13,13).Value = [double]"18.134913"
$ws.Cells.Item(13,14).Value = [double]"54.404739"
$ws.Cells.Item(13,15).Value = [double]"0.629990206705752"
$ws.Cells.Item(13,16).Value = [double]"0.6299902067057519"
$ws.Cells.Item(13,17).Value = [double]"1.552070484134"
$ws.Cells.Item(13,18).Value = [double]"13.968634357206"
$ws.Cells.Item(13,19).Value = [double]"0.0009743714257576927"
$ws.Cells.Item(13,20).Value = [double]"0.0009743714257576925"

$ws.Cells.Item(14,7).Value = [double]"8.835068666666666"
$ws.Cells.Item(14,8).Value = [double]"26.505206"
$ws.Cells.Item(14,9).Value = [double]"0.1596631549860458"
$ws.Cells.Item(14,10).Value = [double]"0.1596631549860458"
$ws.Cells.Item(14,13).Value = [double]"0.043414"
$ws.Cells.Item(14,15).Value = [double]"0.00150816245073376"
$ws.Cells.Item(14,16).Value = [double]"0.00150816245073376"
$ws.Cells.Item(14,17).Value = [double]"0.3835656710946667"
$ws.Cells.Item(14,18).Value = [double]"3.452091039852"
$ws.Cells.Item(14,19).Value = [double]"0.0002407979751156389"
$ws.Cells.Item(14,20).Value = [double]"0.0002407979751156389"

$ws.Cells.Item(15,7).Value = [double]"8.835068666666666"
$ws.Cells.Item(15,8).Value = [double]"26.505206"
$ws.Cells.Item(15,9).Value = [double]"0.1596631549860458"
$ws.Cells.Item(15,10).Value = [double]"0.1596631549860458"
$ws.Cells.Item(15,13).Value = [double]"7.896729"
$ws.Cells.Item(15,14).Value = [double]"23.690187"
$ws.Cells.Item(15,15).Value = [double]"0.2743251062196607"
$ws.Cells.Item(15,16).Value = [double]"0.2743251062196607"
$ws.Cells.Item(15,17).Value = [double]"69.76814295705799"
$ws.Cells.Item(15,18).Value = [double]"627.913286613522"
$ws.Cells.Item(15,19).Value = [double]"0.04379961195091317"
$ws.Cells.Item(15,20).Value = [double]"0.04379961195091316"

$ws.Cells.Item(16,7).Value = [double]"8.835068666666666"
$ws.Cells.Item(16,8).Value = [double]"26.505206"
$ws.Cells.Item(16,9).Value = [double]"0.1596631549860458"
$ws.Cells.Item(16,10).Value = [double]"0.1596631549860458"
$ws.Cells.Item(16,13).Value = [double]"2.710967666666666"
$ws.Cells.Item(16,14).Value = [double]"8.132902999999999"
$ws.Cells.Item(16,15).Value = [double]"0.09417652462385363"
$ws.Cells.Item(16,16).Value = [double]"0.09417652462385362"
$ws.Cells.Item(16,17).Value = [double]"23.95158548811311"
$ws.Cells.Item(16,18).Value = [double]"215.564269393018"
$ws.Cells.Item(16,19).Value = [double]"0.0150365210470655"
$ws.Cells.Item(16,20).Value = [double]"0.0150365210470655"

$ws.Cells.Item(17,7).Value = [double]"8.835068666666666"
$ws.Cells.Item(17,8).Value = [double]"26.505206"
$ws.Cells.Item(17,9).Value = [double]"0.1596631549860458"
$ws.Cells.Item(17,10).Value = [double]"0.1596631549860458"
$ws.Cells.Item(17,13).Value = [double]"18.134913"
$ws.Cells.Item(17,14).Value = [double]"54.404739"
$ws.Cells.Item(17,15).Value = [double]"0.629990206705752"
$ws.Cells.Item(17,16).Value = [double]"0.6299902067057519"
$ws.Cells.Item(17,17).Value = [double]"160.223201619026"
$ws.Cells.Item(17,18).Value = [double]"1442.008814571234"
$ws.Cells.Item(17,19).Value = [double]"0.1005862240129515"
$ws.Cells.Item(17,20).Value = [double]"0.1005862240129515"

$ws.Cells.Item(18,5).Value = [double]"2"
$ws.Cells.Item(18,6).Value = [double]"0.6666666666666666"
$ws.Cells.Item(18,7).Value = [double]"0.3169006666666667"
$ws.Cells.Item(18,8).Value = [double]"0.950702"
$ws.Cells.Item(18,9).Value = [double]"0.00572687798659417"
$ws.Cells.Item(18,10).Value = [double]"0.00572687798659417"
$ws.Cells.Item(18,13).Value = [double]"0.043414"
$ws.Cells.Item(18,15).Value = [double]"0.00150816245073376"
$ws.Cells.Item(18,16).Value = [double]"0.00150816245073376"
$ws.Cells.Item(18,17).Value = [double]"0.01375792554266667"
$ws.Cells.Item(18,18).Value = [double]"0.123821329884"
$ws.Cells.Item(18,19).Value = [double]"8.637062339315082E-06"
$ws.Cells.Item(18,20).Value = [double]"8.637062339315082E-06"

$ws.Cells.Item(19,5).Value = [double]"2"
$ws.Cells.Item(19,6).Value = [double]"0.6666666666666666"
$ws.Cells.Item(19,7).Value = [double]"0.3169006666666667"
$ws.Cells.Item(19,8).Value = [double]"0.950702"
$ws.Cells.Item(19,9).Value = [double]"0.00572687798659417"
$ws.Cells.Item(19,10).Value = [double]"0.00572687798659417"
$ws.Cells.Item(19,13).Value = [double]"7.896729"
$ws.Cells.Item(19,14).Value = [double]"23.690187"
$ws.Cells.Item(19,15).Value = [double]"0.2743251062196607"
$ws.Cells.Item(19,16).Value = [double]"0.2743251062196607"
$ws.Cells.Item(19,17).Value = [double]"2.502478684586"
$ws.Cells.Item(19,18).Value = [double]"22.522308161274"
$ws.Cells.Item(19,19).Value = [double]"0.001571026411979482"
$ws.Cells.Item(19,20).Value = [double]"0.001571026411979482"

$ws.Cells.Item(20,5).Value = [double]"2"
$ws.Cells.Item(20,6).Value = [double]"0.6666666666666666"
$ws.Cells.Item(20,7).Value = [double]"0.3169006666666667"
$ws.Cells.Item(20,8).Value = [double]"0.950702"
$ws.Cells.Item(20,9).Value = [double]"0.00572687798659417"
$ws.Cells.Item(20,10).Value = [double]"0.00572687798659417"
$ws.Cells.Item(20,13).Value = [double]"2.710967666666666"
$ws.Cells.Item(20,14).Value = [double]"8.132902999999999"
$ws.Cells.Item(20,15).Value = [double]"0.09417652462385363"
$ws.Cells.Item(20,16).Value = [double]"0.09417652462385362"
$ws.Cells.Item(20,17).Value = [double]"0.8591074608784444"
$ws.Cells.Item(20,18).Value = [double]"7.731967147905999"
$ws.Cells.Item(20,19).Value = [double]"0.0005393374657222912"
$ws.Cells.Item(20,20).Value = [double]"0.0005393374657222911"

$ws.Cells.Item(21,5).Value = [double]"2"
$ws.Cells.Item(21,6).Value = [double]"0.6666666666666666"
$ws.Cells.Item(21,7).Value = [double]"0.3169006666666667"
$ws.Cells.Item(21,8).Value = [double]"0.950702"
$ws.Cells.Item(21,9).Value = [double]"0.00572687798659417"
$ws.Cells.Item(21,10).Value = [double]"0.00572687798659417"
$ws.Cells.Item(21,13).Value = [double]"18.134913"
$ws.Cells.Item(21,14).Value = [double]"54.404739"
$ws.Cells.Item(21,15).Value = [double]"0.629990206705752"
$ws.Cells.Item(21,16).Value = [double]"0.6299902067057519"
$ws.Cells.Item(21,17).Value = [double]"5.746966019642"
$ws.Cells.Item(21,18).Value = [double]"51.722694176778"
$ws.Cells.Item(21,19).Value = [double]"0.003607877046553082"
$ws.Cells.Item(21,20).Value = [double]"0.003607877046553081"

$ws.Cells.Item(22,5).Value = [double]"3"
$ws.Cells.Item(22,6).Value = [double]"1"
$ws.Cells.Item(22,7).Value = [double]"2.411519"
$ws.Cells.Item(22,8).Value = [double]"7.234557"
$ws.Cells.Item(22,9).Value = [double]"0.04357982335796155"
$ws.Cells.Item(22,10).Value = [double]"0.04357982335796154"
$ws.Cells.Item(22,13).Value = [double]"0.043414"
$ws.Cells.Item(22,15).Value = [double]"0.00150816245073376"
$ws.Cells.Item(22,16).Value = [double]"0.00150816245073376"
$ws.Cells.Item(22,17).Value = [double]"0.104693685866"
$ws.Cells.Item(22,18).Value = [double]"0.9422431727939999"
$ws.Cells.Item(22,19).Value = [double]"6.572545319808762E-05"
$ws.Cells.Item(22,20).Value = [double]"6.572545319808762E-05"

$ws.Cells.Item(23,5).Value = [double]"3"
$ws.Cells.Item(23,6).Value = [double]"1"
$ws.Cells.Item(23,7).Value = [double]"2.411519"
$ws.Cells.Item(23,8).Value = [double]"7.234557"
$ws.Cells.Item(23,9).Value = [double]"0.04357982335796155"
$ws.Cells.Item(23,10).Value = [double]"0.04357982335796154"
$ws.Cells.Item(23,13).Value = [double]"7.896729"
$ws.Cells.Item(23,14).Value = [double]"23.690187"
$ws.Cells.Item(23,15).Value = [double]"0.2743251062196607"
$ws.Cells.Item(23,16).Value = [double]"0.2743251062196607"
$ws.Cells.Item(23,17).Value = [double]"19.043112021351"
$ws.Cells.Item(23,18).Value = [double]"171.388008192159"
$ws.Cells.Item(23,19).Value = [double]"0.01195503967170685"
$ws.Cells.Item(23,20).Value = [double]"0.01195503967170685"

$ws.Cells.Item(24,5).Value = [double]"3"
$ws.Cells.Item(24,6).Value = [double]"1"
$ws.Cells.Item(24,7).Value = [double]"2.411519"
$ws.Cells.Item(24,8).Value = [double]"7.234557"
$ws.Cells.Item(24,9).Value = [double]"0.04357982335796155"
$ws.Cells.Item(24,10).Value = [double]"0.04357982335796154"
$ws.Cells.Item(24,13).Value = [double]"2.710967666666666"
$ws.Cells.Item(24,14).Value = [double]"8.132902999999999"
$ws.Cells.Item(24,15).Value = [double]"0.09417652462385363"
$ws.Cells.Item(24,16).Value = [double]"0.09417652462385362"
$ws.Cells.Item(24,17).Value = [double]"6.537550036552332"
$ws.Cells.Item(24,18).Value = [double]"58.83795032897099"
$ws.Cells.Item(24,19).Value = [double]"0.004104196307574257"
$ws.Cells.Item(24,20).Value = [double]"0.004104196307574256"

$ws.Cells.Item(25,5).Value = [double]"3"
$ws.Cells.Item(25,6).Value = [double]"1"
$ws.Cells.Item(25,7).Value = [double]"2.411519"
$ws.Cells.Item(25,8).Value = [double]"7.234557"
$ws.Cells.Item(25,9).Value = [double]"0.04357982335796155"
$ws.Cells.Item(25,10).Value = [double]"0.04357982335796154"
$ws.Cells.Item(25,13).Value = [double]"18.134913"
$ws.Cells.Item(25,14).Value = [double]"54.404739"
$ws.Cells.Item(25,15).Value = [double]"0.629990206705752"
$ws.Cells.Item(25,16).Value = [double]"0.6299902067057519"
$ws.Cells.Item(25,17).Value = [double]"43.73268726284699"
$ws.Cells.Item(25,18).Value = [double]"393.594185365623"
$ws.Cells.Item(25,19).Value = [double]"0.02745486192548235"
$ws.Cells.Item(25,20).Value = [double]"0.02745486192548234"
